$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15236
$ws1.Range("F4").Value = 81
$ws1.Range("F5").Value = 1579
$ws1.Range("F8").Value = 136
$ws1.Range("G8").Value = 25
$ws1.Range("F10").Value = 8406
$ws1.Range("F11").Value = 995
$ws1.Range("F12").Value = 54
$ws1.Range("F13").Value = 15
$ws1.Range("F15").Value = 1295
$ws1.Range("F16").Value = 62
$ws1.Range("F19").Value = 8740
$ws1.Range("F20").Value = 146
$ws1.Range("F22").Value = 201
$ws1.Range("F23").Value = 161
$ws1.Range("F24").Value = 327
$ws1.Range("F25").Value = 5836
$ws1.Range("F26").Value = 1029
$ws1.Range("F29").Value = 81

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15236
$ws4.Range("F4").Value = 81
$ws4.Range("F5").Value = 1579
$ws4.Range("F9").Value = 136
$ws4.Range("G9").Value = 25
$ws4.Range("F11").Value = 8407
$ws4.Range("F12").Value = 995
$ws4.Range("F13").Value = 54
$ws4.Range("F14").Value = 15
$ws4.Range("F16").Value = 1295
$ws4.Range("F17").Value = 62
$ws4.Range("F22").Value = 8741
$ws4.Range("F23").Value = 146
$ws4.Range("F25").Value = 201
$ws4.Range("F26").Value = 161
$ws4.Range("F27").Value = 327
$ws4.Range("F28").Value = 5836
$ws4.Range("F29").Value = 1029
$ws4.Range("F32").Value = 81
